$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 238.80952
$ws.Range("J9").Value = 659.75
$ws.Range("L9").Value = 659.75
$ws.Range("N9").Value = -997.75
$ws.Range("H51").Value = 10223.5625
$ws.Range("I51").Value = 8121.75
$ws.Range("J51").Value = 10924.167
$ws.Range("K51").Value = 8121.75
$ws.Range("L51").Value = 10924.167
$ws.Range("M51").Value = -7637.75
$ws.Range("N51").Value = -11892.167
$ws.Range("H111").Value = 2559.8
$ws.Range("I111").Value = 2559.8
$ws.Range("K111").Value = 7679.400000000001
$ws.Range("M111").Value = -4612.400000000001
$ws.Range("H127").Value = 172598540
$ws.Range("I127").Value = 90910140
$ws.Range("K127").Value = 272730420
$ws.Range("M127").Value = -272725460
$ws.Range("H138").Value = 2860.2856
$ws.Range("I138").Value = 1088.2307
$ws.Range("J138").Value = 3500.1943
$ws.Range("K138").Value = 3264.6921
$ws.Range("L138").Value = 10500.5829
$ws.Range("M138").Value = 1875.3079
$ws.Range("N138").Value = -20780.5829
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 2508.1428
$ws.Range("I141").Value = 2253
$ws.Range("K141").Value = 6759
$ws.Range("M141").Value = -1579

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3017.4866
$ws.Range("I61").Value = 2119.6785
$ws.Range("K61").Value = 2119.6785
$ws.Range("M61").Value = -1907.6785
$ws.Range("H136").Value = 3017.4866
$ws.Range("I136").Value = 2119.6785
$ws.Range("K136").Value = 6359.0355
$ws.Range("M136").Value = -3809.0355

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2218.5557
$ws.Range("I105").Value = 2049.2
$ws.Range("K105").Value = 2049.2
$ws.Range("M105").Value = -302.1999999999998
$ws.Range("H107").Value = 2076.1482
$ws.Range("I107").Value = 2139.0417
$ws.Range("J107").Value = 1573
$ws.Range("K107").Value = 2139.0417
$ws.Range("L107").Value = 1573
$ws.Range("M107").Value = -219.0417000000002
$ws.Range("N107").Value = -5413
$ws.Range("H134").Value = 5514.5347
$ws.Range("I134").Value = 5249.385
$ws.Range("K134").Value = 15748.155
$ws.Range("M134").Value = -13213.155

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 314781.94
$ws.Range("I31").Value = 478020.53
$ws.Range("J31").Value = 3144.6365
$ws.Range("K31").Value = 478020.53
$ws.Range("L31").Value = 3144.6365
$ws.Range("M31").Value = -477725.53
$ws.Range("N31").Value = -3734.6365
$ws.Range("H34").Value = 314781.94
$ws.Range("I34").Value = 478020.53
$ws.Range("J34").Value = 3144.6365
$ws.Range("K34").Value = 478020.53
$ws.Range("L34").Value = 3144.6365
$ws.Range("M34").Value = -477818.53
$ws.Range("N34").Value = -3548.6365
$ws.Range("H58").Value = 2598.658
$ws.Range("I58").Value = 2220
$ws.Range("K58").Value = 2220
$ws.Range("M58").Value = -2017
$ws.Range("H105").Value = 4539.372
$ws.Range("I105").Value = 1783.2778
$ws.Range("J105").Value = 6523.76
$ws.Range("K105").Value = 1783.2778
$ws.Range("L105").Value = 6523.76
$ws.Range("M105").Value = -36.27780000000007
$ws.Range("N105").Value = -10017.76
$ws.Range("H107").Value = 3474.18
$ws.Range("J107").Value = 5602.241
$ws.Range("L107").Value = 5602.241
$ws.Range("N107").Value = -9442.241
$ws.Range("H122").Value = 3739.6924
$ws.Range("I122").Value = 3812.889
$ws.Range("J122").Value = 3575
$ws.Range("K122").Value = 11438.667
$ws.Range("L122").Value = 10725
$ws.Range("M122").Value = -8988.667000000001
$ws.Range("N122").Value = -15625
$ws.Range("H136").Value = 2598.658
$ws.Range("I136").Value = 2220
$ws.Range("K136").Value = 6660
$ws.Range("M136").Value = -4110

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 214418200
$ws.Range("I76").Value = 300182750
$ws.Range("J76").Value = 6849.5
$ws.Range("K76").Value = 900548250
$ws.Range("L76").Value = 20548.5
$ws.Range("M76").Value = -900547867
$ws.Range("N76").Value = -21314.5
$ws.Range("H79").Value = 214418200
$ws.Range("I79").Value = 300182750
$ws.Range("J79").Value = 6849.5
$ws.Range("K79").Value = 900548250
$ws.Range("L79").Value = 20548.5
$ws.Range("M79").Value = -900546924
$ws.Range("N79").Value = -23200.5
$ws.Range("H81").Value = 499.5
$ws.Range("I81").Value = 499.5
$ws.Range("K81").Value = 1498.5
$ws.Range("M81").Value = -375.5
$ws.Range("H84").Value = 499.5
$ws.Range("I84").Value = 499.5
$ws.Range("K84").Value = 4495.5
$ws.Range("M84").Value = 1120.5
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 65710.19
$ws.Range("I102").Value = 3424.2
$ws.Range("K102").Value = 3424.2
$ws.Range("M102").Value = -1802.2
$ws.Range("H132").Value = 34699.53
$ws.Range("I132").Value = 39438.82
$ws.Range("K132").Value = 118316.46
$ws.Range("M132").Value = -115786.46
$ws.Range("H136").Value = 7000
$ws.Range("J136").Value = 7000
$ws.Range("L136").Value = 21000
$ws.Range("N136").Value = -26100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3126
$ws.Range("I122").Value = 3347
$ws.Range("K122").Value = 10041
$ws.Range("M122").Value = -7591
$ws.Range("H132").Value = 2811.25
$ws.Range("I132").Value = 2658.9614
$ws.Range("J132").Value = 4791
$ws.Range("K132").Value = 7976.8842
$ws.Range("L132").Value = 14373
$ws.Range("M132").Value = -5446.8842
$ws.Range("N132").Value = -19433

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 49999.668
$ws.Range("J41").Value = 49999.668
$ws.Range("L41").Value = 49999.668
$ws.Range("N41").Value = -50779.668
$ws.Range("H122").Value = 4220.923
$ws.Range("I122").Value = 4279.273
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 12837.819
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -10387.819
$ws.Range("N122").Value = -16600
$ws.Range("H126").Value = 4470.909
$ws.Range("I126").Value = 4240.4287
$ws.Range("K126").Value = 12721.2861
$ws.Range("M126").Value = -10251.2861
$ws.Range("H132").Value = 1228.375
$ws.Range("I132").Value = 1217
$ws.Range("K132").Value = 3651
$ws.Range("M132").Value = -1121
$ws.Range("H136").Value = 589349.3
$ws.Range("I136").Value = 589349.3
$ws.Range("K136").Value = 1768047.9
$ws.Range("M136").Value = -1765497.9
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
